$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

$ws = $sheets.Item(1)
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H32").Value = 10560.267
$ws.Range("I32").Value = 9188.799999999999
$ws.Range("J32").Value = 11246
$ws.Range("K32").Value = 9188.799999999999
$ws.Range("L32").Value = 11246
$ws.Range("M32").Value = -8862.799999999999
$ws.Range("N32").Value = -11898
$ws.Range("H62").Value = 5466.1665
$ws.Range("I62").Value = 3718.8
$ws.Range("J62").Value = 6714.2856
$ws.Range("K62").Value = 3718.8
$ws.Range("L62").Value = 6714.2856
$ws.Range("M62").Value = -3094.8
$ws.Range("N62").Value = -7962.2856
$ws.Range("H65").Value = 5466.1665
$ws.Range("I65").Value = 3718.8
$ws.Range("J65").Value = 6714.2856
$ws.Range("K65").Value = 18594
$ws.Range("L65").Value = 33571.428
$ws.Range("M65").Value = -15474
$ws.Range("N65").Value = -39811.428
$ws.Range("H68").Value = 68723
$ws.Range("J68").Value = 68723
$ws.Range("L68").Value = 68723
$ws.Range("N68").Value = -70221
$ws.Range("H71").Value = 68723
$ws.Range("J71").Value = 68723
$ws.Range("L71").Value = 206169
$ws.Range("N71").Value = -213657
$ws.Range("H80").Value = 983.5
$ws.Range("I80").Value = 751
$ws.Range("K80").Value = 2253
$ws.Range("M80").Value = -1255
$ws.Range("H83").Value = 983.5
$ws.Range("I83").Value = 751
$ws.Range("K83").Value = 6759
$ws.Range("M83").Value = -1767
$ws.Range("H96").Value = 1710.75
$ws.Range("I96").Value = 1211.1
$ws.Range("K96").Value = 3633.3
$ws.Range("M96").Value = -2260.3
$ws.Range("H98").Value = 58825800
$ws.Range("I98").Value = 62502236
$ws.Range("K98").Value = 62502236
$ws.Range("M98").Value = -62500738
$ws.Range("H113").Value = 125002220
$ws.Range("I113").Value = 50002000
$ws.Range("J113").Value = 200002450
$ws.Range("K113").Value = 50002000
$ws.Range("L113").Value = 200002450
$ws.Range("M113").Value = -49998746
$ws.Range("N113").Value = -200008958
$ws.Range("H122").Value = 58825800
$ws.Range("I122").Value = 62502236
$ws.Range("K122").Value = 187506708
$ws.Range("M122").Value = -187504258
$ws.Range("H132").Value = 1735.1333
$ws.Range("I132").Value = 1648.1786
$ws.Range("K132").Value = 4944.5358
$ws.Range("M132").Value = -2414.5358
$ws.Range("H137").Value = 4964.174
$ws.Range("J137").Value = 9258.333000000001
$ws.Range("L137").Value = 27774.999
$ws.Range("N137").Value = -32874.999
$ws.Range("H141").Value = 2247.5
$ws.Range("I141").Value = 1996.6666
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 5989.9998
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -809.9997999999996
$ws.Range("N141").Value = -19360
$ws = $sheets.Item(2)
$ws.Range("H45").Value = 2410.4614
$ws.Range("I45").Value = 2011
$ws.Range("K45").Value = 2011
$ws.Range("M45").Value = -1634
$ws.Range("H102").Value = 20758.285
$ws.Range("I102").Value = 20758.285
$ws.Range("K102").Value = 20758.285
$ws.Range("M102").Value = -19136.285
$ws.Range("H122").Value = 3128.6155
$ws.Range("I122").Value = 1959
$ws.Range("K122").Value = 5877
$ws.Range("M122").Value = -3427
$ws.Range("H132").Value = 8096.1665
$ws.Range("I132").Value = 2081.3333
$ws.Range("J132").Value = 14111
$ws.Range("K132").Value = 6243.999899999999
$ws.Range("L132").Value = 42333
$ws.Range("M132").Value = -3713.999899999999
$ws.Range("N132").Value = -47393
$ws = $sheets.Item(3)
$ws.Range("H16").Value = 1599
$ws.Range("I16").Value = 1599
$ws.Range("K16").Value = 1599
$ws.Range("M16").Value = -1429
$ws.Range("H86").Value = 1074.8182
$ws.Range("I86").Value = 1030.762
$ws.Range("K86").Value = 1030.762
$ws.Range("M86").Value = 92.23800000000006
$ws.Range("H89").Value = 1074.8182
$ws.Range("I89").Value = 1030.762
$ws.Range("K89").Value = 5153.809999999999
$ws.Range("M89").Value = 462.1900000000005
$ws.Range("H99").Value = 2270
$ws.Range("I99").Value = 1538.4546
$ws.Range("K99").Value = 1538.4546
$ws.Range("M99").Value = -40.45460000000003
$ws.Range("H134").Value = 29698.842
$ws.Range("I134").Value = 3125.8484
$ws.Range("J134").Value = 205080.6
$ws.Range("K134").Value = 9377.5452
$ws.Range("L134").Value = 615241.8
$ws.Range("M134").Value = -6842.5452
$ws.Range("N134").Value = -620311.8
$ws = $sheets.Item(4)
$ws.Range("H7").Value = 1457.625
$ws.Range("I7").Value = 110.833336
$ws.Range("K7").Value = 110.833336
$ws.Range("M7").Value = 2.166663999999997
$ws.Range("H16").Value = 2029.2222
$ws.Range("J16").Value = 2199.6
$ws.Range("L16").Value = 2199.6
$ws.Range("N16").Value = -2773.6
$ws.Range("H22").Value = 476.7143
$ws.Range("I22").Value = 277.6
$ws.Range("J22").Value = 974.5
$ws.Range("K22").Value = 277.6
$ws.Range("L22").Value = 974.5
$ws.Range("M22").Value = 72.39999999999998
$ws.Range("N22").Value = -1674.5
$ws.Range("H31").Value = 1556808.5
$ws.Range("I31").Value = 2050.6667
$ws.Range("K31").Value = 2050.6667
$ws.Range("M31").Value = -1755.6667
$ws.Range("H34").Value = 1556808.5
$ws.Range("I34").Value = 2050.6667
$ws.Range("K34").Value = 2050.6667
$ws.Range("M34").Value = -1848.6667
$ws.Range("H58").Value = 1433.9375
$ws.Range("I58").Value = 1566.0834
$ws.Range("J58").Value = 1037.5
$ws.Range("K58").Value = 1566.0834
$ws.Range("L58").Value = 1037.5
$ws.Range("M58").Value = -1363.0834
$ws.Range("N58").Value = -1443.5
$ws.Range("H113").Value = 2029.2222
$ws.Range("J113").Value = 2199.6
$ws.Range("L113").Value = 2199.6
$ws.Range("N113").Value = -6539.6
$ws.Range("H132").Value = 2036.3572
$ws.Range("I132").Value = 2116.077
$ws.Range("K132").Value = 6348.231000000001
$ws.Range("M132").Value = -3818.231000000001
$ws.Range("H134").Value = 336089.6
$ws.Range("I134").Value = 418087.38
$ws.Range("K134").Value = 1254262.14
$ws.Range("M134").Value = -1251727.14
$ws.Range("H136").Value = 1433.9375
$ws.Range("I136").Value = 1566.0834
$ws.Range("J136").Value = 1037.5
$ws.Range("K136").Value = 4698.2502
$ws.Range("L136").Value = 3112.5
$ws.Range("M136").Value = -2148.2502
$ws.Range("N136").Value = -8212.5
$ws = $sheets.Item(5)
$ws.Range("H137").Value = 7689.615
$ws.Range("J137").Value = 10076.625
$ws.Range("L137").Value = 30229.875
$ws.Range("N137").Value = -40429.875
$ws = $sheets.Item(6)
$ws.Range("H122").Value = 20666.334
$ws.Range("I122").Value = 17749.75
$ws.Range("J122").Value = 26499.5
$ws.Range("K122").Value = 53249.25
$ws.Range("L122").Value = 79498.5
$ws.Range("M122").Value = -50799.25
$ws.Range("N122").Value = -84398.5
$ws.Range("H132").Value = 71431310
$ws.Range("I132").Value = 76925820
$ws.Range("J132").Value = 2689
$ws.Range("K132").Value = 230777460
$ws.Range("L132").Value = 8067
$ws.Range("M132").Value = -230774930
$ws.Range("N132").Value = -13127
$ws = $sheets.Item(7)
$ws.Range("H7").Value = 51642.855
$ws.Range("I7").Value = 3173.5
$ws.Range("J7").Value = 116268.664
$ws.Range("K7").Value = 3173.5
$ws.Range("L7").Value = 116268.664
$ws.Range("M7").Value = -3061.5
$ws.Range("N7").Value = -116492.664
$ws.Range("H22").Value = 1030.3334
$ws.Range("I22").Value = 1034.7
$ws.Range("J22").Value = 1021.6
$ws.Range("K22").Value = 1034.7
$ws.Range("L22").Value = 1021.6
$ws.Range("M22").Value = -739.7
$ws.Range("N22").Value = -1611.6
$ws.Range("H27").Value = 1030.3334
$ws.Range("I27").Value = 1034.7
$ws.Range("J27").Value = 1021.6
$ws.Range("K27").Value = 1034.7
$ws.Range("L27").Value = 1021.6
$ws.Range("M27").Value = -927.7
$ws.Range("N27").Value = -1235.6
$ws.Range("H93").Value = 71429760
$ws.Range("I93").Value = 90910240
$ws.Range("K93").Value = 90910240
$ws.Range("M93").Value = -90908992
$ws.Range("H100").Value = 3077.4707
$ws.Range("I100").Value = 3087.8
$ws.Range("K100").Value = 3087.8
$ws.Range("M100").Value = -2546.8
$ws.Range("H122").Value = 5556.2173
$ws.Range("I122").Value = 4576.1763
$ws.Range("J122").Value = 8333
$ws.Range("K122").Value = 13728.5289
$ws.Range("L122").Value = 24999
$ws.Range("M122").Value = -11278.5289
$ws.Range("N122").Value = -29899
$ws.Range("H126").Value = 51642.855
$ws.Range("I126").Value = 3173.5
$ws.Range("J126").Value = 116268.664
$ws.Range("K126").Value = 9520.5
$ws.Range("L126").Value = 348805.992
$ws.Range("M126").Value = -7050.5
$ws.Range("N126").Value = -353745.992
$ws.Range("H132").Value = 136783.94
$ws.Range("I132").Value = 75125.28999999999
$ws.Range("J132").Value = 1000005
$ws.Range("K132").Value = 225375.87
$ws.Range("L132").Value = 3000015
$ws.Range("M132").Value = -222845.87
$ws.Range("N132").Value = -3005075
$ws.Range("H136").Value = 156864.92
$ws.Range("I136").Value = 169350.67
$ws.Range("K136").Value = 508052.01
$ws.Range("M136").Value = -505502.01
$ws = $sheets.Item(8)
$ws.Range("H81").Value = 1001
$ws.Range("I81").Value = 1000
$ws.Range("K81").Value = 2000
$ws.Range("M81").Value = -939
$ws.Range("H84").Value = 1001
$ws.Range("I84").Value = 1000
$ws.Range("K84").Value = 10000
$ws.Range("M84").Value = -4696
$ws.Range("H126").Value = 7474.84
$ws.Range("I126").Value = 6547.85
$ws.Range("K126").Value = 19643.55
$ws.Range("M126").Value = -17173.55

Write-Output "done"